# "added number into genSymbolVals" - the generated symbol list gained a new
# entry, which shifted the iteration/shuffle order of the existing rows.
# The underlying data set (rows 2-25) is unchanged as a set, only the row
# order differs. Re-write the sheet with the new row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(201,  9, 30, 15, 45, 30),
    @(1202, 2, 10, 10, 10, 10),
    @(601,  9, 60, 67, 60, 42),
    @(1201, 2, 10, 10, 10, 10),
    @(1203, 3, 15, 15, 15, 15),
    @(701,  3, 90, 45, 97, 15),
    @(801,  3, 67, 65, 52, 45),
    @(101,  9, 30, 15, 60, 15),
    @(301,  6, 45, 30, 60, 45),
    @(501,  9, 52, 30, 75, 45),
    @(401,  9, 48, 67, 75, 45),
    @(901, 16, 15, 45, 60, 60),
    @(902,  1,  0,  0,  0,  0),
    @(1001,18, 30, 75, 60, 72),
    @(2,    0,  2,  2,  2,  2),
    @(1,    0,  2,  2,  2,  2),
    @(502,  0,  4,  0,  0,  0),
    @(802,  0,  4,  5,  4,  0),
    @(3,    0,  3,  3,  3,  3),
    @(1101, 0, 15, 30, 30,  0),
    @(602,  0,  0,  4,  0,  9),
    @(402,  0,  0,  4,  0,  0),
    @(702,  0,  0,  0,  4,  0),
    @(1002, 0,  0,  0,  0,  9)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
